$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.520.72'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.726.71'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +12.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.09'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '657.32'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.444'
$ws.Range("E8").Value = '  +4.72%  '
$ws.Range("E9").Value = '  +3.96%  '
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = '3.725.33'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000315'
$ws.Range("E12").Value = '  +16.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.77'
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.85'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '4.423.20'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = '97.353.37'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.86'
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("D19").Value = '3.718.08'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.10'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.88'
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.540'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '528.93'
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.48'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("E25").Value = '  +9.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '117.14'
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("E28").Value = '  +26.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.45'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.81'
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.04'
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  +2.92%  '
$ws.Range("E34").Value = '  -2.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.14'
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.599'
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '639.93'
$ws.Range("E38").Value = '  -3.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.77'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.169'
$ws.Range("E41").Value = '  +4.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.499'
$ws.Range("E42").Value = '  +11.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.84'
$ws.Range("E43").Value = '  -4.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.92'
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.02'
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.968'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0458'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  +2.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.79'
$ws.Range("E49").Value = '  +2.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.65'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.36'
$ws.Range("E51").Value = '  +4.03%  '
